$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to 8e8a11af-c2bc-445c-a82c-9a46df9ff85d.md
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 corresponds to 8e8a11af-c2bc-445c-a82c-9a46df9ff85d.md
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-01-25 06:32:33"

# de-de sheet: row 3 corresponds to 8e8a11af-c2bc-445c-a82c-9a46df9ff85d.md
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-01-25 06:32:52"
